$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 15000
$ws.Range("J7").Value = 15000
$ws.Range("L7").Value = 15000
$ws.Range("N7").Value = -15224
$ws.Range("H14").Value = 15000
$ws.Range("J14").Value = 15000
$ws.Range("L14").Value = 15000
$ws.Range("N14").Value = -15382
$ws.Range("H28").Value = 482.23077
$ws.Range("J28").Value = 562.375
$ws.Range("L28").Value = 562.375
$ws.Range("N28").Value = -1532.375
$ws.Range("H51").Value = 16294.263
$ws.Range("I51").Value = 7915.1665
$ws.Range("J51").Value = 20161.54
$ws.Range("K51").Value = 7915.1665
$ws.Range("L51").Value = 20161.54
$ws.Range("M51").Value = -7431.1665
$ws.Range("N51").Value = -21129.54
$ws.Range("H113").Value = 79930.92
$ws.Range("I113").Value = 2812.125
$ws.Range("J113").Value = 203321
$ws.Range("K113").Value = 2812.125
$ws.Range("L113").Value = 203321
$ws.Range("M113").Value = 441.875
$ws.Range("N113").Value = -209829
$ws.Range("H132").Value = 1124.7826
$ws.Range("I132").Value = 1146.1904
$ws.Range("K132").Value = 3438.5712
$ws.Range("M132").Value = -908.5711999999999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5753597.5
$ws.Range("I2").Value = 7078120
$ws.Range("K2").Value = 7078120
$ws.Range("M2").Value = -7078007
$ws.Range("H45").Value = 89999.5
$ws.Range("I45").Value = 89999.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 89999.5
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -89622.5
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 6955.25
$ws.Range("I61").Value = 5794
$ws.Range("J61").Value = 9510
$ws.Range("K61").Value = 5794
$ws.Range("L61").Value = 9510
$ws.Range("M61").Value = -5582
$ws.Range("N61").Value = -9934
$ws.Range("H96").Value = 32472
$ws.Range("J96").Value = 32472
$ws.Range("L96").Value = 32472
$ws.Range("N96").Value = -37964
$ws.Range("H116").Value = 5753597.5
$ws.Range("I116").Value = 7078120
$ws.Range("K116").Value = 7078120
$ws.Range("M116").Value = -7075826
$ws.Range("H132").Value = 7306.0835
$ws.Range("I132").Value = 4349.5386
$ws.Range("J132").Value = 10800.182
$ws.Range("K132").Value = 13048.6158
$ws.Range("L132").Value = 32400.546
$ws.Range("M132").Value = -10518.6158
$ws.Range("N132").Value = -37460.546
$ws.Range("H136").Value = 6955.25
$ws.Range("I136").Value = 5794
$ws.Range("J136").Value = 9510
$ws.Range("K136").Value = 17382
$ws.Range("L136").Value = 28530
$ws.Range("M136").Value = -14832
$ws.Range("N136").Value = -33630

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5753597.5
$ws.Range("I3").Value = 7078120
$ws.Range("K3").Value = 7078120
$ws.Range("M3").Value = -7078006
$ws.Range("H134").Value = 5206.1113
$ws.Range("I134").Value = 3014
$ws.Range("K134").Value = 9042
$ws.Range("M134").Value = -6507

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8029.636
$ws.Range("I16").Value = 5842.3335
$ws.Range("K16").Value = 5842.3335
$ws.Range("M16").Value = -5555.3335
$ws.Range("H25").Value = 4066.6667
$ws.Range("I25").Value = 4100
$ws.Range("K25").Value = 4100
$ws.Range("M25").Value = -3926
$ws.Range("H31").Value = 37041310
$ws.Range("I31").Value = 142858190
$ws.Range("J31").Value = 5404.5
$ws.Range("K31").Value = 142858190
$ws.Range("L31").Value = 5404.5
$ws.Range("M31").Value = -142857895
$ws.Range("N31").Value = -5994.5
$ws.Range("H34").Value = 37041310
$ws.Range("I34").Value = 142858190
$ws.Range("J34").Value = 5404.5
$ws.Range("K34").Value = 142858190
$ws.Range("L34").Value = 5404.5
$ws.Range("M34").Value = -142857988
$ws.Range("N34").Value = -5808.5
$ws.Range("H58").Value = 5202.5356
$ws.Range("I58").Value = 3360.8572
$ws.Range("J58").Value = 7044.2144
$ws.Range("K58").Value = 3360.8572
$ws.Range("L58").Value = 7044.2144
$ws.Range("M58").Value = -3157.8572
$ws.Range("N58").Value = -7450.2144
$ws.Range("H59").Value = 40654
$ws.Range("I59").Value = 31276
$ws.Range("J59").Value = 46012.855
$ws.Range("K59").Value = 31276
$ws.Range("L59").Value = 46012.855
$ws.Range("M59").Value = -30131
$ws.Range("N59").Value = -48302.855
$ws.Range("H60").Value = 49773
$ws.Range("I60").Value = 49546
$ws.Range("J60").Value = 50000
$ws.Range("K60").Value = 49546
$ws.Range("L60").Value = 50000
$ws.Range("M60").Value = -49035
$ws.Range("N60").Value = -51022
$ws.Range("H99").Value = 5259.25
$ws.Range("I99").Value = 4837.778
$ws.Range("J99").Value = 5801.143
$ws.Range("K99").Value = 4837.778
$ws.Range("L99").Value = 5801.143
$ws.Range("M99").Value = -3339.778
$ws.Range("N99").Value = -8797.143
$ws.Range("H113").Value = 8029.636
$ws.Range("I113").Value = 5842.3335
$ws.Range("K113").Value = 5842.3335
$ws.Range("M113").Value = -3672.3335
$ws.Range("H126").Value = 5259.25
$ws.Range("I126").Value = 4837.778
$ws.Range("J126").Value = 5801.143
$ws.Range("K126").Value = 14513.334
$ws.Range("L126").Value = 17403.429
$ws.Range("M126").Value = -12043.334
$ws.Range("N126").Value = -22343.429
$ws.Range("H134").Value = 6296.1377
$ws.Range("I134").Value = 6566.696
$ws.Range("J134").Value = 5259
$ws.Range("K134").Value = 19700.088
$ws.Range("L134").Value = 15777
$ws.Range("M134").Value = -17165.088
$ws.Range("N134").Value = -20847
$ws.Range("H136").Value = 5202.5356
$ws.Range("I136").Value = 3360.8572
$ws.Range("J136").Value = 7044.2144
$ws.Range("K136").Value = 10082.5716
$ws.Range("L136").Value = 21132.6432
$ws.Range("M136").Value = -7532.571599999999
$ws.Range("N136").Value = -26232.6432

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 42240996
$ws.Range("I4").Value = 75357470
$ws.Range("K4").Value = 226072410
$ws.Range("M4").Value = -226072298
$ws.Range("H74").Value = 10714
$ws.Range("J74").Value = 10939.8
$ws.Range("L74").Value = 32819.39999999999
$ws.Range("N74").Value = -34941.39999999999
$ws.Range("H77").Value = 10714
$ws.Range("J77").Value = 10939.8
$ws.Range("L77").Value = 98458.2
$ws.Range("N77").Value = -109066.2
$ws.Range("H130").Value = 6452.5
$ws.Range("I130").Value = 6452.5
$ws.Range("K130").Value = 19357.5
$ws.Range("M130").Value = -14337.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1197.8182
$ws.Range("I102").Value = 1018.35
$ws.Range("K102").Value = 1018.35
$ws.Range("M102").Value = 603.65
$ws.Range("H113").Value = 52270.715
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4013.5
$ws.Range("I40").Value = 2444
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 2444
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = -2308
$ws.Range("N40").Value = -15272
$ws.Range("H55").Value = 708.38464
$ws.Range("I55").Value = 466.875
$ws.Range("J55").Value = 1094.8
$ws.Range("K55").Value = 466.875
$ws.Range("L55").Value = 1094.8
$ws.Range("M55").Value = -293.875
$ws.Range("N55").Value = -1440.8
$ws.Range("H61").Value = 6761.143
$ws.Range("I61").Value = 6400.55
$ws.Range("K61").Value = 6400.55
$ws.Range("M61").Value = -6198.55
$ws.Range("H82").Value = 2988.3076
$ws.Range("I82").Value = 1968
$ws.Range("K82").Value = 1968
$ws.Range("M82").Value = -1607
$ws.Range("H85").Value = 2988.3076
$ws.Range("I85").Value = 1968
$ws.Range("K85").Value = 1968
$ws.Range("M85").Value = -720
$ws.Range("H113").Value = 6761.143
$ws.Range("I113").Value = 6400.55
$ws.Range("K113").Value = 6400.55
$ws.Range("M113").Value = -4230.55
$ws.Range("H132").Value = 6944.4375
$ws.Range("I132").Value = 6463.727
$ws.Range("K132").Value = 19391.181
$ws.Range("M132").Value = -16861.181
$ws.Range("H136").Value = 4453.607
$ws.Range("I136").Value = 2702.3125
$ws.Range("J136").Value = 6788.6665
$ws.Range("K136").Value = 8106.9375
$ws.Range("L136").Value = 20365.9995
$ws.Range("M136").Value = -5556.9375
$ws.Range("N136").Value = -25465.9995

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 174251
$ws.Range("I2").Value = 11200.75
$ws.Range("J2").Value = 500351.5
$ws.Range("K2").Value = 11200.75
$ws.Range("L2").Value = 500351.5
$ws.Range("M2").Value = -11088.75
$ws.Range("N2").Value = -500575.5
$ws.Range("H113").Value = 597.35297
$ws.Range("I113").Value = 682.1
$ws.Range("J113").Value = 476.2857
$ws.Range("K113").Value = 2046.3
$ws.Range("L113").Value = 1428.8571
$ws.Range("M113").Value = 123.6999999999998
$ws.Range("N113").Value = -5768.8571
$ws.Range("H122").Value = 4814.8
$ws.Range("I122").Value = 5075.524
$ws.Range("J122").Value = 3446
$ws.Range("K122").Value = 15226.572
$ws.Range("L122").Value = 10338
$ws.Range("M122").Value = -12776.572
$ws.Range("N122").Value = -15238
$ws.Range("H126").Value = 4615.05
$ws.Range("I126").Value = 4615.05
$ws.Range("K126").Value = 13845.15
$ws.Range("M126").Value = -11375.15
